$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant
$xlPasteFormats = -4122

# --- Row 3: add 2020 / 2021 / 2022 headers in K3:M3 (same style as J3) ---
$ws.Range("J3").Copy()
$ws.Range("K3:M3").PasteSpecial($xlPasteFormats)
$ws.Range("K3").Value = 2020
$ws.Range("L3").Value = 2021
$ws.Range("M3").Value = 2022

# --- Row 4: add K4/L4/M4 values (same style as J4) ---
$ws.Range("J4").Copy()
$ws.Range("K4:M4").PasteSpecial($xlPasteFormats)
$ws.Range("K4").Value = 308
$ws.Range("L4").Value = 212.1
$ws.Range("M4").Value = 723.8

# --- Row 5: E5:L5 become "-" placeholders with a new right-aligned style; M5 gets the
#     same style but stays empty. Base the new style on D4's (fontId 3, vertical center)
#     then switch horizontal alignment to right. ---
$ws.Range("D4").Copy()
$ws.Range("E5:M5").PasteSpecial($xlPasteFormats)
$ws.Range("E5:M5").HorizontalAlignment = -4152
$ws.Range("E5").Value = "-"
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = "-"
$ws.Range("H5").Value = "-"
$ws.Range("I5").Value = "-"
$ws.Range("J5").Value = "-"
$ws.Range("K5").Value = "-"
$ws.Range("L5").Value = "-"

# --- Row 6: add J6/K6/L6 values (same style as I6); M6 stays empty with that style ---
$ws.Range("I6").Copy()
$ws.Range("J6:M6").PasteSpecial($xlPasteFormats)
$ws.Range("J6").Value = 9.8000000000000007
$ws.Range("K6").Value = 9.8000000000000007
$ws.Range("L6").Value = 9.8000000000000007

# --- Row 7: J7 (already styled) gets a value; K7:M7 added with the same style as I7 ---
$ws.Range("I7").Copy()
$ws.Range("K7:M7").PasteSpecial($xlPasteFormats)
$ws.Range("J7").Value = 64
$ws.Range("K7").Value = 64
$ws.Range("L7").Value = 64
$ws.Range("M7").Value = 64

# --- Update the selected cell to match the saved view state ---
$ws.Range("M14").Select()

Write-Output "done"
